$d = $word.ActiveDocument

# 1. Title paragraph: "第2节 DIY图形化控制小车" -> "第1节 DIY图形化控制小车".
#    In the original file this text was split across three runs (with a
#    hidden "_GoBack" bookmark sitting between the "2" and the "节..." run).
#    Replacing the whole phrase collapses it back down to a single run,
#    matching a normal "select & retype" edit, and drops the stale
#    "_GoBack" bookmark that lived in the middle of it.
$d.Content.Find.Execute("第2节 DIY图形化控制小车", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "第1节 DIY图形化控制小车", 2) | Out-Null

# 2. Word keeps exactly one "_GoBack" bookmark around, marking the location
#    of the most recent edit. Since the last edit of this session was made
#    in the (empty) "适用年级" value cell of the info table, re-create the
#    bookmark there - this is what pushes the OLE_LINK73/74/75/76 bookmark
#    ids back down/up by one.
$table = $d.Tables.Item(1)
$cell = $table.Cell(2, 2)
$cell.Range.Bookmarks.Add("_GoBack") | Out-Null
